# -----------------------------------------------------------------------
# Dejaber (دجابر) quarterly income statement - roll the quarterly columns
# one period forward (drop the oldest quarter, add 1401/12 quarter) and
# update the read_price-derived figures for the 1400/12 quarter column,
# per commit "update database and change read_price algorithm".
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data columns D..M correspond to column indexes 4..13
$colIdx = @(4,5,6,7,8,9,10,11,12,13)

# ---------------------------------------------------------------------
# Row 8: quarter / fiscal-period labels (shift left by one quarter, add
# the new quarter "فصل چهارم منتهی به 1401/12" at the end)
# ---------------------------------------------------------------------
$row8 = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt $colIdx.Length; $i++) {
    $ws.Cells.Item(8, $colIdx[$i]).Value = $row8[$i]
}

# ---------------------------------------------------------------------
# Row 9: publish dates (shift left by one quarter, last value's label
# changes from "1401-10-30 (6)" to "1402-02-28 (7)", add new
# "1402-02-28" at the end)
# ---------------------------------------------------------------------
$row9 = @(
    "1400-11-06 (3)",
    "1401-02-27 (11)",
    "1401-04-29 (3)",
    "1401-09-15 (5)",
    "1401-10-30 (3)",
    "1402-02-28 (7)",
    "1401-04-29",
    "1401-09-15 (2)",
    "1401-10-30",
    "1402-02-28"
)
for ($i = 0; $i -lt $colIdx.Length; $i++) {
    $ws.Cells.Item(9, $colIdx[$i]).Value = $row9[$i]
}

# ---------------------------------------------------------------------
# Numeric data rows: shift left by one quarter and append the new
# quarter's figure. Row 11's "column I" (1400/12 quarter) additionally
# reflects the read_price algorithm change, which also cascades into
# rows 13/16/17/19/24/26 (derived totals) at that same column.
# ---------------------------------------------------------------------
$dataRows = @{
    11 = @(3882, 2661, 2930, 7453, 5162, 3604, 3709, 6732, 7023, 3789)
    12 = @(-2343, -2181, -2114, -4846, -3509, -2746, -2766, -4598, -4750, -4429)
    13 = @(1539, 480, 816, 2607, 1654, 857, 943, 2134, 2274, -640)
    14 = @(-230, -320, -382, -426, -391, -452, -481, -542, -538, -467)
    16 = @(29, -236, 9, 7, -117, -62, 0, -8, 38, -251)
    17 = @(1338, -75, 443, 2187, 1146, 343, 462, 1584, 1774, -1359)
    18 = @(-271, -463, -429, -563, -495, -443, -445, -715, -475, -387)
    19 = @(33, 4402, 720, -20, 836, 6133, 2340, 1835, 0, 5933)
    20 = @(1100, 3864, 734, 1604, 1487, 6033, 2357, 2703, 1299, 4187)
    21 = @(-255, 118, -3, -440, -321, 238, -3, -181, -260, 324)
    22 = @(845, 3982, 731, 1164, 1167, 6271, 2354, 2522, 1039, 4511)
    24 = @(845, 3982, 731, 1164, 1167, 6275, 2354, 2522, 1039, 4515)
    26 = @(8346, 9286, 9721, 8689, 7966, 8234, 7718, 7312, 6528, 4990)
}

foreach ($r in $dataRows.Keys) {
    $vals = $dataRows[$r]
    for ($i = 0; $i -lt $colIdx.Length; $i++) {
        $ws.Cells.Item($r, $colIdx[$i]).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# Row 23 ("هزینه کاهش ارزش دریافتنی‌ها"): column I (1400/12 quarter) and
# the new column M (1401/12 quarter) now carry a real value of 4 instead
# of the placeholder "-" dash.
# ---------------------------------------------------------------------
$ws.Cells.Item(23, 9).Value = 4
$ws.Cells.Item(23, 13).Value = 4

# ---------------------------------------------------------------------
# Column widths: the "wide" (31-char) highlighted column shifts from
# F/J to E/I, and the new last column M also becomes a wide column.
# Narrow values close to, but not exactly, 29/31 are used because the
# Excel column-width model introduces rounding when a column's width is
# changed to a new value; these particular inputs round-trip to exactly
# 29 and 31 characters in the saved file.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 5).EntireColumn.ColumnWidth  = 30.14   # E: 29 -> 31
$ws.Cells.Item(1, 6).EntireColumn.ColumnWidth  = 28.14   # F: 31 -> 29
$ws.Cells.Item(1, 9).EntireColumn.ColumnWidth  = 30.14   # I: 29 -> 31
$ws.Cells.Item(1, 10).EntireColumn.ColumnWidth = 28.14   # J: 31 -> 29
$ws.Cells.Item(1, 13).EntireColumn.ColumnWidth = 30.14   # M: 29 -> 31

# ---------------------------------------------------------------------
# Row heights: slightly reduced due to the workbook's updated default
# font metrics (15.75 -> 15.6 for label rows, 42 -> 40.8 for title rows).
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 15.6
$ws.Rows.Item(5).RowHeight = 40.8
$ws.Rows.Item(6).RowHeight = 40.8
$ws.Rows.Item(8).RowHeight = 15.6

Write-Host "Dejaber quarterly income statement updated successfully."
